# Auto-generated edit script applying stock report updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 146
$ws.Range("G6").Value = 4362.48
$ws.Range("F9").Value = 97
$ws.Range("G9").Value = 2868.29
$ws.Range("B10").Value = 38238.21
$ws.Range("F71").Value = 19
$ws.Range("G71").Value = 4765.01
$ws.Range("F86").Value = 3
$ws.Range("G86").Value = 338.52
$ws.Range("F90").Value = 76
$ws.Range("G90").Value = 10256.2
$ws.Range("F109").Value = 111
$ws.Range("G109").Value = 13927.17
$ws.Range("B114").Value = 243050.24
$ws.Range("F143").Value = 6
$ws.Range("G143").Value = 583.2
$ws.Range("B152").Value = 21102.29
$ws.Range("F193").Value = 289
$ws.Range("G193").Value = 18727.2
$ws.Range("F197").Value = 60
$ws.Range("G197").Value = 2789.4
$ws.Range("F198").Value = 15
$ws.Range("G198").Value = 1127.1
$ws.Range("B200").Value = 43906.85
$ws.Range("F237").Value = 31
$ws.Range("G237").Value = 1549.07
$ws.Range("B246").Value = 64973
$ws.Range("E246").Value = 35.4
$ws.Range("F246").Value = 17
$ws.Range("G246").Value = 566.1
$ws.Range("B247").Value = 48706
$ws.Range("E247").Value = 39.8
$ws.Range("F247").Value = -144
$ws.Range("G247").Value = -4795.2
$ws.Range("F272").Value = 66
$ws.Range("G272").Value = 5788.2
$ws.Range("F273").Value = 9
$ws.Range("G273").Value = 3808.89
$ws.Range("B274").Value = 71805.43
$ws.Range("F285").Value = 2
$ws.Range("G285").Value = 263.66
$ws.Range("B315").Value = 60325
$ws.Range("E315").Value = 151.57
$ws.Range("F315").Value = -102
$ws.Range("G315").Value = -12939.72
$ws.Range("B316").Value = 63560
$ws.Range("E316").Value = 134.87
$ws.Range("F316").Value = 1
$ws.Range("G316").Value = 126.86
$ws.Range("F325").Value = 168
$ws.Range("G325").Value = 23153.76
$ws.Range("F328").Value = 411
$ws.Range("G328").Value = 8643.33
$ws.Range("B339").Value = 276419.59
$ws.Range("F355").Value = 78
$ws.Range("G355").Value = 2514.72
$ws.Range("B361").Value = 12039.02
$ws.Range("F363").Value = 24
$ws.Range("G363").Value = 510
$ws.Range("F379").Value = 82
$ws.Range("G379").Value = 7625.18
$ws.Range("B395").Value = 233773.12
$ws.Range("F426").Value = 57
$ws.Range("G426").Value = 5506.2
$ws.Range("B430").Value = 41404.99
$ws.Range("F433").Value = 45
$ws.Range("G433").Value = 2245.05
$ws.Range("B448").Value = 37897.65
$ws.Range("B472").Value = 45695
$ws.Range("E472").Value = 23.58
$ws.Range("F472").Value = -36
$ws.Range("G472").Value = -710.28
$ws.Range("B473").Value = 64915
$ws.Range("E473").Value = 20.98
$ws.Range("F473").Value = 0
$ws.Range("G473").Value = 0
$ws.Range("F479").Value = 134
$ws.Range("G479").Value = 2173.48
$ws.Range("F486").Value = 151
$ws.Range("G486").Value = 1985.65
$ws.Range("B490").Value = 65067
$ws.Range("E490").Value = 15.65
$ws.Range("F490").Value = 215
$ws.Range("G490").Value = 3166.95
$ws.Range("B491").Value = 53595
$ws.Range("E491").Value = 17.61
$ws.Range("F491").Value = -335
$ws.Range("G491").Value = -4934.55
$ws.Range("B492").Value = -12138.67
$ws.Range("F498").Value = 101
$ws.Range("G498").Value = 6226.65
$ws.Range("B508").Value = 7990.3
$ws.Range("F579").Value = 9
$ws.Range("G579").Value = 1725.12
$ws.Range("F580").Value = 49
$ws.Range("G580").Value = 1298.99
$ws.Range("B587").Value = 19873.19
$ws.Range("B596").Value = 64836
$ws.Range("E596").Value = 104.71
$ws.Range("F596").Value = 0
$ws.Range("G596").Value = 0
$ws.Range("B597").Value = 60031
$ws.Range("E597").Value = 111.69
$ws.Range("F597").Value = -5
$ws.Range("G597").Value = -492.5
$ws.Range("F655").Value = 318
$ws.Range("G655").Value = 25560.84
$ws.Range("B656").Value = 33961.6
$ws.Range("F703").Value = 4
$ws.Range("G703").Value = 326.24
$ws.Range("F706").Value = 122
$ws.Range("G706").Value = 17461.86
$ws.Range("F707").Value = 13
$ws.Range("G707").Value = 1060.28
$ws.Range("F710").Value = 36
$ws.Range("G710").Value = 2724.48
$ws.Range("F712").Value = 36
$ws.Range("G712").Value = 781.92
$ws.Range("F717").Value = 5
$ws.Range("G717").Value = 554.15
$ws.Range("F718").Value = 245
$ws.Range("G718").Value = 33077.45
$ws.Range("F720").Value = 164
$ws.Range("G720").Value = 19796.44
$ws.Range("B721").Value = 93407.76
$ws.Range("F725").Value = 20
$ws.Range("G725").Value = 3277.8
$ws.Range("F729").Value = 12
$ws.Range("G729").Value = 1804.92
$ws.Range("B737").Value = 65362
$ws.Range("F737").Value = 44
$ws.Range("G737").Value = 1798.28
$ws.Range("B738").Value = 65079
$ws.Range("F738").Value = 21
$ws.Range("G738").Value = 858.27
$ws.Range("F740").Value = 62
$ws.Range("G740").Value = 2049.72
$ws.Range("F743").Value = 127
$ws.Range("G743").Value = 6336.03
$ws.Range("F746").Value = 90
$ws.Range("G746").Value = 21762
$ws.Range("B748").Value = 64077.33
$ws.Range("F773").Value = 2852
$ws.Range("G773").Value = 465189.72
$ws.Range("F774").Value = 30
$ws.Range("G774").Value = 5285.1
$ws.Range("F775").Value = 558
$ws.Range("G775").Value = 157841.46
$ws.Range("F776").Value = 446
$ws.Range("G776").Value = 64513.9
$ws.Range("F779").Value = 96
$ws.Range("G779").Value = 12343.68
$ws.Range("B780").Value = 705808.53
$ws.Range("F783").Value = 66
$ws.Range("G783").Value = 9636.66
$ws.Range("F790").Value = 190
$ws.Range("G790").Value = 30240.4
$ws.Range("B797").Value = 69326.71
$ws.Range("B798").Value = 2661679.4
$ws.Range("B799").Value = 2661679.4
